$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "A thought particle can be object particle denoted with V ..." ->
#    "A thought particle can be an object particle denoted with V ..."
#    with "object particle" italicised.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Start = 0
$rng.End = 0
$ok = $rng.Find.Execute("an be object particle denoted with V", $true, $false, $false, $false, $false, $true, 1, $false, "an be an object particle denoted with V", 2)
if (-not $ok) { throw "step1: phrase not found" }

$rng = $d.Content
$rng.Start = 0
$rng.End = 0
$ok = $rng.Find.Execute("object particle", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) { throw "step2: phrase not found" }
$rng.Font.Italic = $true

# ---------------------------------------------------------------------------
# 2) "... or connecting particle denoted with A ..." ->
#    "... or a connecting particle denoted with A ..."
#    with "connecting particle" italicised.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Start = 0
$rng.End = 0
$ok = $rng.Find.Execute("or connecting particle denoted with A", $true, $false, $false, $false, $false, $true, 1, $false, "or a connecting particle denoted with A", 2)
if (-not $ok) { throw "step3: phrase not found" }

$rng = $d.Content
$rng.Start = 0
$rng.End = 0
$ok = $rng.Find.Execute("connecting particle", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) { throw "step4: phrase not found" }
$rng.Font.Italic = $true

# ---------------------------------------------------------------------------
# 3) Italicise "V-particle" and "A-particle" (each unique in the document).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Start = 0
$rng.End = 0
$ok = $rng.Find.Execute("V-particle", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) { throw "step5: phrase not found" }
$rng.Font.Italic = $true

$rng = $d.Content
$rng.Start = 0
$rng.End = 0
$ok = $rng.Find.Execute("A-particle", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) { throw "step6: phrase not found" }
$rng.Font.Italic = $true

# ---------------------------------------------------------------------------
# 4) Add the new closing paragraph right after the "Every connecting
#    particle signature encodes ..." paragraph.
# ---------------------------------------------------------------------------
$inserted = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Every connecting particle signature encodes*") {
        $p.Range.InsertParagraphAfter()
        $newPara = $d.Paragraphs($i + 1)
        $newPara.Range.Text = "Obviously the signature of every compound thought particle is mxn matrix."
        $inserted = $true
        break
    }
}
if (-not $inserted) { throw "step7: anchor paragraph not found" }
